# linjeimpedanser.xlsx - update shunt admittance/impedance columns with
# correct computed values and add a new "Shunt Admittans (p.u.) / 2" column.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Header row (row 1) ----
$ws.Range("F1").Value = "Shunt Impedans (ohm)"
$ws.Range("G1").Value = "Shunt Admittans (S)"
$ws.Range("H1").Value = "Shunt Admittans (p.u.)"
$ws.Range("I1").Value = "Shunt Admittans (p.u.) / 2"

# Give the new header cell (I1) the same style as the other header cells.
$ws.Range("H1").Copy() | Out-Null
$ws.Range("I1").PasteSpecial(-4122) | Out-Null  # xlPasteFormats
$ws.Range("I1").Value = "Shunt Admittans (p.u.) / 2"
$excel.CutCopyMode = $false

# ---- Data rows (row 2 - row 10) ----
$ws.Range("F2").Value = "-15977.0058j"
$ws.Range("G2").Value = "0.0001j"
$ws.Range("H2").Value = "0.0056j"
$ws.Range("I2").Value = "0.0028j"

$ws.Range("F3").Value = "-19461.3528j"
$ws.Range("G3").Value = "0.0001j"
$ws.Range("H3").Value = "0.0046j"
$ws.Range("I3").Value = "0.0023j"

$ws.Range("F4").Value = "-7654.2559j"
$ws.Range("G4").Value = "0.0001j"
$ws.Range("H4").Value = "0.0118j"
$ws.Range("I4").Value = "0.0059j"

$ws.Range("F5").Value = "-23280.1789j"
$ws.Range("G5").Value = "0j"
$ws.Range("H5").Value = "0.0039j"
$ws.Range("I5").Value = "0.0019j"

$ws.Range("F6").Value = "-3921.618j"
$ws.Range("G6").Value = "0.0003j"
$ws.Range("H6").Value = "0.0229j"
$ws.Range("I6").Value = "0.0115j"

$ws.Range("F7").Value = "-5132.7059j"
$ws.Range("G7").Value = "0.0002j"
$ws.Range("H7").Value = "0.0175j"
$ws.Range("I7").Value = "0.0088j"

$ws.Range("F8").Value = "-7406.3448j"
$ws.Range("G8").Value = "0.0001j"
$ws.Range("H8").Value = "0.0122j"
$ws.Range("I8").Value = "0.0061j"

$ws.Range("F9").Value = "-3458.161j"
$ws.Range("G9").Value = "0.0003j"
$ws.Range("H9").Value = "0.026j"
$ws.Range("I9").Value = "0.013j"

$ws.Range("F10").Value = "-3135.1622j"
$ws.Range("G10").Value = "0.0003j"
$ws.Range("H10").Value = "0.0287j"
$ws.Range("I10").Value = "0.0144j"

# ---- Column widths (best effort; engine quantizes to 1/6 character units) ----
$ws.Columns.Item(2).ColumnWidth = 15.333333333333334
$ws.Columns.Item(3).ColumnWidth = 16.666666666666668
$ws.Columns.Item(4).ColumnWidth = 17.166666666666668
$ws.Columns.Item(5).ColumnWidth = 13.833333333333334
$ws.Columns.Item(6).ColumnWidth = 20.5
$ws.Columns.Item(7).ColumnWidth = 17.666666666666668
$ws.Columns.Item(8).ColumnWidth = 18.5
$ws.Columns.Item(9).ColumnWidth = 22.833333333333332

# ---- Selection, matching the workbook state at save time ----
$ws.Range("F14").Select() | Out-Null
